$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert 5 new rows above the existing content, pushing everything down.
$ws.Rows("1:5").Insert()

# 2. New header/description row (was inserted blank row 2): merged B2:E2
#    with the new descriptive text, left-aligned + wrap text.
$descCell = $ws.Range("B2")
$descCell.HorizontalAlignment = -4131
$descCell.WrapText = $true
$ws.Range("B2:E2").Merge()
$ws.Range("B2").Value = "This excel file demonstrates how a `$10,000 principal accumulates at different compound interest rates over time."
$ws.Rows(2).RowHeight = 43

# 3. Row 3 is hidden, row 4 gets a small custom height (spacer rows).
$ws.Rows(3).Hidden = $true
$ws.Rows(4).RowHeight = 11.5

# 4. Update the chart's series source ranges to point at the shifted data
#    (rows 4-11 -> rows 9-16), keeping the cached plotted values untouched.
$chartObj = $ws.ChartObjects().Item(1)
$chart = $chartObj.Chart
$sc = $chart.SeriesCollection()
for ($i = 1; $i -le $sc.Count(); $i++) {
  $s = $sc.Item($i)
  $f = $s.Formula()
  $newF = $f.Replace("`$4:", "`$9:").Replace("`$11,", "`$16,").Replace("`$11)", "`$16)")
  $s.Formula = $newF
}

# 5. Move the chart down by the height of the 5 inserted rows (default row
#    height 14.5pt x 5 rows = 72.5pt) so it keeps its position relative to
#    the data table underneath it.
$chartObj.Top = $chartObj.Top() + 72.5

# 6. Restore the active selection to the new description cell.
$ws.Range("B2:E2").Select()
